$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (it sat right after "CN_08_01_CO").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Replace the "Enunciado" sentence and re-insert the "_GoBack" bookmark
#    right after the new sentence (before the trailing ". ").
$d.Content.Find.Execute(
    "Coloca la palabra correspondiente al sitio señalado en la imagen que se muestra a continuación. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Relaciona cada palabra con su ubicación en la imagen. ", 2) | Out-Null

# Re-create the _GoBack bookmark right after "...en la imagen" (before the ". ").
$r = $d.Content
$r.Find.Execute("Relaciona cada palabra con su ubicación en la imagen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bm = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

# 3. Remove the "No hay ninguna palabra que se repita en la imagen." sentence,
#    leaving its paragraph empty, then merge away the following blank paragraph
#    (the blank paragraph's own mark is deleted, so the sentence's paragraph
#    -- now empty -- is the one that survives).
$target = $d.Content
$target.Find.Execute(
    "No hay ninguna palabra que se repita en la imagen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null

$sentStart = $target.Start
$sentEnd = $target.End

# Delete the paragraph mark that ends the following blank paragraph, merging
# it into the sentence's paragraph.
$mergeRange = $d.Range($sentEnd + 1, $sentEnd + 2)
$mergeRange.Delete()

# Now delete the sentence text itself, leaving a single empty paragraph.
$sentRange = $d.Range($sentStart, $sentEnd)
$sentRange.Delete()

# 4. Flip the "Mostrar al inicio del ejercicio ... (S/N)" answer from S to N.
#    Locate the whole sentence first so the match is unique (there is another
#    "(S/N) N" elsewhere), then replace only the trailing " S" run in place so
#    the "(S/N)" run's own formatting/boundaries stay untouched.
$sentence = $d.Content
$sentence.Find.Execute(
    "Mostrar al inicio del ejercicio ventana Más información (S/N) S",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$tailRange = $d.Range($sentence.End - 2, $sentence.End)
$tailRange.Text = " N"
